$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Metadata / content corrections (per commit message: "2013 download files metadata corrected") ---

# F1: was a HYPERLINK() formula whose display text was the old page title;
# now a plain text label referencing the new shared string.
$ws.Range("F1").Value = "Metadata - Single European Sky Portal"

# F1's (hyperlink-styled) font: color + face updated.
$ws.Range("F1").Font.Color = 13391121   # RGB(0x11,0x55,0xCC) -> BGR-packed OLE color
$ws.Range("F1").Font.Name = "Arial"

# F2: contact e-mail address text corrected.
$ws.Range("F2").Value = "pru-support@eurocontrol.int"

# --- Outline defaults now explicitly written (sheetPr/outlinePr) ---
$outline = $ws.Outline
$outline.SummaryRow = 0
$outline.SummaryColumn = 0

# --- Column width adjustments ---
$ws.Columns.Item(1).ColumnWidth = 12.250666666666667
$ws.Columns.Item(2).ColumnWidth = 14.250666666666667
$ws.Columns.Item(3).ColumnWidth = 8.750666666666666
$ws.Columns.Item(4).ColumnWidth = 13.250666666666667
$ws.Columns.Item(5).ColumnWidth = 11.083666666666666
$ws.Columns.Item(6).ColumnWidth = 11.250666666666667

# --- Restore frozen panes / selection state (unaffected by the content edits above) ---
$ws.Range("B6").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true | Out-Null
